$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.731.40'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.898.21'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.07%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9950'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.58%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.88'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('B6').Value = 'XRP'
$ws.Range('C6').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5424'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +14.83%  '
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9967'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2932'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06569'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.69'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '100.45'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.91%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7651'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.25%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07846'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.884.82'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.64%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.278'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '285.35'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.618.94'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.26'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007566'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9994'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.127.42'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.378'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9952'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.65%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.470'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.221'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.87'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.927'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1013'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.342'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.508'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.47%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.281'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.83%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.216'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04885'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.142'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7055'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.781'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01918'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.877'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.339'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '76.09'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.990'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.74%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4293'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8445'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9986'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.992'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.49'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.090'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.37%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.40'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05798'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3989'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.79%  '
